# The commit swaps the presentation's design theme: the deck's active
# design ("Integral", backed by the theme part the slide master/
# presentation point at) is replaced by the stock "Office Theme" colour
# palette (the palette that used to live only behind the notes master).
#
# PowerPoint exposes the active theme's colours through
# Master.ColorScheme (an indexed 1..12 collection: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). Re-pointing every one of those swatches to
# the "Office Theme" values reproduces the same palette switch an author
# gets by picking "Office Theme" from the Design gallery.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = $r + ($g * 256) + ($b * 65536)
    $scheme.Colors($index).RGB = $bgr
}

# Index -> (theme element, target "Office Theme" RGB)
Set-ThemeColor $colorScheme 1  "000000"   # dk1
Set-ThemeColor $colorScheme 2  "FFFFFF"   # lt1
Set-ThemeColor $colorScheme 3  "44546A"   # dk2
Set-ThemeColor $colorScheme 4  "E7E6E6"   # lt2
Set-ThemeColor $colorScheme 5  "5B9BD5"   # accent1
Set-ThemeColor $colorScheme 6  "ED7D31"   # accent2
Set-ThemeColor $colorScheme 7  "A5A5A5"   # accent3
Set-ThemeColor $colorScheme 8  "FFC000"   # accent4
Set-ThemeColor $colorScheme 9  "4472C4"   # accent5
Set-ThemeColor $colorScheme 10 "70AD47"   # accent6
Set-ThemeColor $colorScheme 11 "0563C1"   # hlink
Set-ThemeColor $colorScheme 12 "954F72"   # folHlink
